$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title placeholder: "Here" " " "is" " " "a" " " "single" " " "header" -> "Here is a single header"
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "X"
$titleRange.Text = "Here is a single header"

# Speaker notes body: "and" " " "here" " " "are" " " "some" " " "notes" -> "and here are some notes"
$notesRange = $s.NotesPage.Shapes.Item(2).TextFrame.TextRange
$notesRange.Text = "and here are some notes"
